$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 42.48700833333334
$ws.Range("H2").Value2 = 127.461025
$ws.Range("I2").Value2 = 0.8741865936964721
$ws.Range("J2").Value2 = 0.877455058515614
$ws.Range("M2").Value2 = 0.8151449999999999
$ws.Range("N2").Value2 = 2.445435
$ws.Range("O2").Value2 = 0.1271069095499719
$ws.Range("P2").Value2 = 0.1371035811308388
$ws.Range("Q2").Value2 = 34.633072407875
$ws.Range("R2").Value2 = 311.697651670875
$ws.Range("S2").Value2 = 0.1111151562947755
$ws.Range("T2").Value2 = 0.1203022308038604
$ws.Range("G3").Value2 = 42.48700833333334
$ws.Range("H3").Value2 = 127.461025
$ws.Range("I3").Value2 = 0.8741865936964721
$ws.Range("J3").Value2 = 0.877455058515614
$ws.Range("O3").Value2 = 0.4802730342501803
$ws.Range("P3").Value2 = 0.5180454245123947
$ws.Range("Q3").Value2 = 130.8609487055222
$ws.Range("R3").Value2 = 1177.7485383497
$ws.Range("S3").Value2 = 0.4198482478554342
$ws.Range("T3").Value2 = 0.4545615782792694
$ws.Range("G4").Value2 = 42.48700833333334
$ws.Range("H4").Value2 = 127.461025
$ws.Range("I4").Value2 = 0.8741865936964721
$ws.Range("J4").Value2 = 0.877455058515614
$ws.Range("M4").Value2 = 0.5185940000000001
$ws.Range("N4").Value2 = 1.555782
$ws.Range("O4").Value2 = 0.08086522109705406
$ws.Range("P4").Value2 = 0.08722508823947427
$ws.Range("Q4").Value2 = 22.03350759961667
$ws.Range("R4").Value2 = 198.30156839655
$ws.Range("S4").Value2 = 0.07069129217934578
$ws.Range("T4").Value2 = 0.07653609490519749
$ws.Range("G5").Value2 = 42.48700833333334
$ws.Range("H5").Value2 = 127.461025
$ws.Range("I5").Value2 = 0.8741865936964721
$ws.Range("J5").Value2 = 0.877455058515614
$ws.Range("M5").Value2 = 1.402793
$ws.Range("N5").Value2 = 2.805586
$ws.Range("O5").Value2 = 0.2187398352051889
$ws.Range("P5").Value2 = 0.1572954863942594
$ws.Range("Q5").Value2 = 59.60047788094167
$ws.Range("R5").Value2 = 357.60286728565
$ws.Range("S5").Value2 = 0.1912194314437518
$ws.Range("T5").Value2 = 0.1380197202183168
$ws.Range("G6").Value2 = 42.48700833333334
$ws.Range("H6").Value2 = 127.461025
$ws.Range("I6").Value2 = 0.8741865936964721
$ws.Range("J6").Value2 = 0.877455058515614
$ws.Range("M6").Value2 = 0.5965113333333333
$ws.Range("N6").Value2 = 1.789534
$ws.Range("O6").Value2 = 0.09301499989760488
$ws.Range("P6").Value2 = 0.1003304197230327
$ws.Range("Q6").Value2 = 25.34398199026111
$ws.Range("R6").Value2 = 228.09583791235
$ws.Range("S6").Value2 = 0.08131246592316491
$ws.Range("T6").Value2 = 0.0880354343089698
$ws.Range("I7").Value2 = 0.003351874396568939
$ws.Range("J7").Value2 = 0.003364406599215795
$ws.Range("M7").Value2 = 0.8151449999999999
$ws.Range("N7").Value2 = 2.445435
$ws.Range("O7").Value2 = 0.1271069095499719
$ws.Range("P7").Value2 = 0.1371035811308388
$ws.Range("Q7").Value2 = 0.132792826515
$ws.Range("R7").Value2 = 1.195135438635
$ws.Range("S7").Value2 = 0.0004260463957475547
$ws.Range("T7").Value2 = 0.0004612721931327123
$ws.Range("I8").Value2 = 0.003351874396568939
$ws.Range("J8").Value2 = 0.003364406599215795
$ws.Range("O8").Value2 = 0.4802730342501803
$ws.Range("P8").Value2 = 0.5180454245123947
$ws.Range("S8").Value2 = 0.001609814886865657
$ws.Range("T8").Value2 = 0.001742915444923048
$ws.Range("I9").Value2 = 0.003351874396568939
$ws.Range("J9").Value2 = 0.003364406599215795
$ws.Range("M9").Value2 = 0.5185940000000001
$ws.Range("N9").Value2 = 1.555782
$ws.Range("O9").Value2 = 0.08086522109705406
$ws.Range("P9").Value2 = 0.08722508823947427
$ws.Range("Q9").Value2 = 0.08448259275800002
$ws.Range("R9").Value2 = 0.7603433348220001
$ws.Range("S9").Value2 = 0.0002710500641681019
$ws.Range("T9").Value2 = 0.0002934606624900672
$ws.Range("I10").Value2 = 0.003351874396568939
$ws.Range("J10").Value2 = 0.003364406599215795
$ws.Range("M10").Value2 = 1.402793
$ws.Range("N10").Value2 = 2.805586
$ws.Range("O10").Value2 = 0.2187398352051889
$ws.Range("P10").Value2 = 0.1572954863942594
$ws.Range("Q10").Value2 = 0.228524799251
$ws.Range("R10").Value2 = 1.371148795506
$ws.Range("S10").Value2 = 0.0007331884531339819
$ws.Range("T10").Value2 = 0.0005292059724517044
$ws.Range("I11").Value2 = 0.003351874396568939
$ws.Range("J11").Value2 = 0.003364406599215795
$ws.Range("M11").Value2 = 0.5965113333333333
$ws.Range("N11").Value2 = 1.789534
$ws.Range("O11").Value2 = 0.09301499989760488
$ws.Range("P11").Value2 = 0.1003304197230327
$ws.Range("Q11").Value2 = 0.09717587177933332
$ws.Range("R11").Value2 = 0.874582846014
$ws.Range("S11").Value2 = 0.0003117745966536443
$ws.Range("T11").Value2 = 0.0003375523262182618
$ws.Range("G12").Value2 = 3.784599666666667
$ws.Range("H12").Value2 = 11.353799
$ws.Range("I12").Value2 = 0.07786959875243754
$ws.Range("J12").Value2 = 0.07816074259499733
$ws.Range("M12").Value2 = 0.8151449999999999
$ws.Range("N12").Value2 = 2.445435
$ws.Range("O12").Value2 = 0.1271069095499719
$ws.Range("P12").Value2 = 0.1371035811308388
$ws.Range("Q12").Value2 = 3.084997495285
$ws.Range("R12").Value2 = 27.764977457565
$ws.Range("S12").Value2 = 0.009897764045318681
$ws.Range("T12").Value2 = 0.01071611771361983
$ws.Range("G13").Value2 = 3.784599666666667
$ws.Range("H13").Value2 = 11.353799
$ws.Range("I13").Value2 = 0.07786959875243754
$ws.Range("J13").Value2 = 0.07816074259499733
$ws.Range("O13").Value2 = 0.4802730342501803
$ws.Range("P13").Value2 = 0.5180454245123947
$ws.Range("Q13").Value2 = 11.65665275759245
$ws.Range("R13").Value2 = 104.909874818332
$ws.Range("S13").Value2 = 0.03739866846867723
$ws.Range("T13").Value2 = 0.0404908150778294
$ws.Range("G14").Value2 = 3.784599666666667
$ws.Range("H14").Value2 = 11.353799
$ws.Range("I14").Value2 = 0.07786959875243754
$ws.Range("J14").Value2 = 0.07816074259499733
$ws.Range("M14").Value2 = 0.5185940000000001
$ws.Range("N14").Value2 = 1.555782
$ws.Range("O14").Value2 = 0.08086522109705406
$ws.Range("P14").Value2 = 0.08722508823947427
$ws.Range("Q14").Value2 = 1.962670679535334
$ws.Range("R14").Value2 = 17.664036115818
$ws.Range("S14").Value2 = 0.006296942319854747
$ws.Range("T14").Value2 = 0.006817577669711477
$ws.Range("G15").Value2 = 3.784599666666667
$ws.Range("H15").Value2 = 11.353799
$ws.Range("I15").Value2 = 0.07786959875243754
$ws.Range("J15").Value2 = 0.07816074259499733
$ws.Range("M15").Value2 = 1.402793
$ws.Range("N15").Value2 = 2.805586
$ws.Range("O15").Value2 = 0.2187398352051889
$ws.Range("P15").Value2 = 0.1572954863942594
$ws.Range("Q15").Value2 = 5.309009920202334
$ws.Range("R15").Value2 = 31.854059521214
$ws.Range("S15").Value2 = 0.01703318319860237
$ws.Range("T15").Value2 = 0.01229433202341661
$ws.Range("G16").Value2 = 3.784599666666667
$ws.Range("H16").Value2 = 11.353799
$ws.Range("I16").Value2 = 0.07786959875243754
$ws.Range("J16").Value2 = 0.07816074259499733
$ws.Range("M16").Value2 = 0.5965113333333333
$ws.Range("N16").Value2 = 1.789534
$ws.Range("O16").Value2 = 0.09301499989760488
$ws.Range("P16").Value2 = 0.1003304197230327
$ws.Range("Q16").Value2 = 2.257556593296222
$ws.Range("R16").Value2 = 20.318009339666
$ws.Range("S16").Value2 = 0.007243040719984511
$ws.Range("T16").Value2 = 0.007841900110420005
$ws.Range("G17").Value2 = 0.5431155000000001
$ws.Range("H17").Value2 = 1.086231
$ws.Range("I17").Value2 = 0.01117481101996684
$ws.Range("J17").Value2 = 0.007477728079359741
$ws.Range("M17").Value2 = 0.8151449999999999
$ws.Range("N17").Value2 = 2.445435
$ws.Range("O17").Value2 = 0.1271069095499719
$ws.Range("P17").Value2 = 0.1371035811308388
$ws.Range("Q17").Value2 = 0.4427178842475
$ws.Range("R17").Value2 = 2.656307305485
$ws.Range("S17").Value2 = 0.001420395693552954
$ws.Range("T17").Value2 = 0.00102522329840285
$ws.Range("G18").Value2 = 0.5431155000000001
$ws.Range("H18").Value2 = 1.086231
$ws.Range("I18").Value2 = 0.01117481101996684
$ws.Range("J18").Value2 = 0.007477728079359741
$ws.Range("O18").Value2 = 0.4802730342501803
$ws.Range("P18").Value2 = 0.5180454245123947
$ws.Range("Q18").Value2 = 1.672808050618
$ws.Range("R18").Value2 = 10.036848303708
$ws.Range("S18").Value2 = 0.005366960395731827
$ws.Range("T18").Value2 = 0.003873802817260171
$ws.Range("G19").Value2 = 0.5431155000000001
$ws.Range("H19").Value2 = 1.086231
$ws.Range("I19").Value2 = 0.01117481101996684
$ws.Range("J19").Value2 = 0.007477728079359741
$ws.Range("M19").Value2 = 0.5185940000000001
$ws.Range("N19").Value2 = 1.555782
$ws.Range("O19").Value2 = 0.08086522109705406
$ws.Range("P19").Value2 = 0.08722508823947427
$ws.Range("Q19").Value2 = 0.2816564396070001
$ws.Range("R19").Value2 = 1.689938637642
$ws.Range("S19").Value2 = 0.0009036535638474148
$ws.Range("T19").Value2 = 0.0006522454915529478
$ws.Range("G20").Value2 = 0.5431155000000001
$ws.Range("H20").Value2 = 1.086231
$ws.Range("I20").Value2 = 0.01117481101996684
$ws.Range("J20").Value2 = 0.007477728079359741
$ws.Range("M20").Value2 = 1.402793
$ws.Range("N20").Value2 = 2.805586
$ws.Range("O20").Value2 = 0.2187398352051889
$ws.Range("P20").Value2 = 0.1572954863942594
$ws.Range("Q20").Value2 = 0.7618786215915001
$ws.Range("R20").Value2 = 3.047514486366
$ws.Range("S20").Value2 = 0.002444376320956676
$ws.Range("T20").Value2 = 0.001176212875366901
$ws.Range("G21").Value2 = 0.5431155000000001
$ws.Range("H21").Value2 = 1.086231
$ws.Range("I21").Value2 = 0.01117481101996684
$ws.Range("J21").Value2 = 0.007477728079359741
$ws.Range("M21").Value2 = 0.5965113333333333
$ws.Range("N21").Value2 = 1.789534
$ws.Range("O21").Value2 = 0.09301499989760488
$ws.Range("P21").Value2 = 0.1003304197230327
$ws.Range("Q21").Value2 = 0.323974551059
$ws.Range("R21").Value2 = 1.943847306354
$ws.Range("S21").Value2 = 0.001039425045877969
$ws.Range("T21").Value2 = 0.0007502435967768703
$ws.Range("G22").Value2 = 1.624131
$ws.Range("H22").Value2 = 4.872393
$ws.Range("I22").Value2 = 0.03341712213455474
$ws.Range("J22").Value2 = 0.03354206421081321
$ws.Range("M22").Value2 = 0.8151449999999999
$ws.Range("N22").Value2 = 2.445435
$ws.Range("O22").Value2 = 0.1271069095499719
$ws.Range("P22").Value2 = 0.1371035811308388
$ws.Range("Q22").Value2 = 1.323902263995
$ws.Range("R22").Value2 = 11.915120375955
$ws.Range("S22").Value2 = 0.004247547120577211
$ws.Range("T22").Value2 = 0.004598737121823035
$ws.Range("G23").Value2 = 1.624131
$ws.Range("H23").Value2 = 4.872393
$ws.Range("I23").Value2 = 0.03341712213455474
$ws.Range("J23").Value2 = 0.03354206421081321
$ws.Range("O23").Value2 = 0.4802730342501803
$ws.Range("P23").Value2 = 0.5180454245123947
$ws.Range("Q23").Value2 = 5.002360293636
$ws.Range("R23").Value2 = 45.021242642724
$ws.Range("S23").Value2 = 0.01604934264347146
$ws.Range("T23").Value2 = 0.01737631289311273
$ws.Range("G24").Value2 = 1.624131
$ws.Range("H24").Value2 = 4.872393
$ws.Range("I24").Value2 = 0.03341712213455474
$ws.Range("J24").Value2 = 0.03354206421081321
$ws.Range("M24").Value2 = 0.5185940000000001
$ws.Range("N24").Value2 = 1.555782
$ws.Range("O24").Value2 = 0.08086522109705406
$ws.Range("P24").Value2 = 0.08722508823947427
$ws.Range("Q24").Value2 = 0.8422645918140002
$ws.Range("R24").Value2 = 7.580381326326001
$ws.Range("S24").Value2 = 0.002702282969838028
$ws.Range("T24").Value2 = 0.002925709510522294
$ws.Range("G25").Value2 = 1.624131
$ws.Range("H25").Value2 = 4.872393
$ws.Range("I25").Value2 = 0.03341712213455474
$ws.Range("J25").Value2 = 0.03354206421081321
$ws.Range("M25").Value2 = 1.402793
$ws.Range("N25").Value2 = 2.805586
$ws.Range("O25").Value2 = 0.2187398352051889
$ws.Range("P25").Value2 = 0.1572954863942594
$ws.Range("Q25").Value2 = 2.278319597883
$ws.Range("R25").Value2 = 13.669917587298
$ws.Range("S25").Value2 = 0.007309655788744174
$ws.Range("T25").Value2 = 0.005276015304707343
$ws.Range("G26").Value2 = 1.624131
$ws.Range("H26").Value2 = 4.872393
$ws.Range("I26").Value2 = 0.03341712213455474
$ws.Range("J26").Value2 = 0.03354206421081321
$ws.Range("M26").Value2 = 0.5965113333333333
$ws.Range("N26").Value2 = 1.789534
$ws.Range("O26").Value2 = 0.09301499989760488
$ws.Range("P26").Value2 = 0.1003304197230327
$ws.Range("Q26").Value2 = 0.9688125483179999
$ws.Range("R26").Value2 = 8.719312934862
$ws.Range("S26").Value2 = 0.003108293611923859
$ws.Range("T26").Value2 = 0.003365289380647804
